$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Math row (5): status was "Incomplete" -> now "Complete 0.3.1.10b", with line-count 486.
# Copy the visual format used by the other "Complete x.x.x.xb" status cells (e.g. F3)
# onto F5 before changing its text, then fill in the line-count in G5.
$ws.Range("F3").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("F5").Value = "Complete 0.3.1.10b"
$ws.Range("G5").Value = 486

# Dice row (10): status was "In Progress" -> now "Complete 0.3.1.10b", with line-count 386.
$ws.Range("F3").Copy()
$ws.Range("F10").PasteSpecial(-4122)
$ws.Range("F10").Value = "Complete 0.3.1.10b"
$ws.Range("G10").Value = 386

# Xdev row (22): line count was the text "N/A" -> now a real number, 452.
$ws.Range("G22").Value = 452

$excel.CutCopyMode = 0

# Move the active selection cursor (cosmetic, matches author's last position).
$ws.Range("F5").Select()
